# Update data: 9 October 2020
# Adds the September-2020 (date serial 44075) observation to the "Canada"
# sheet and the ten corresponding province rows to the "Province" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "Canada": one new national row -------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

$wsCanada.Range("A10").Value = 44075
$wsCanada.Range("A10").NumberFormat = "d-mmm-yy"
$wsCanada.Range("B10").Value = "Canada"
$wsCanada.Range("B10").NumberFormat = "d-mmm-yy"
$wsCanada.Range("C10").Value = 65.3
$wsCanada.Range("D10").Value = 1832.6

# --- Sheet "Province": ten new province rows ----------------------------
$wsProvince = $wb.Worksheets.Item("Province")

$newRows = @(
    @{ Row = 82; Province = "Newfoundland & Labrador"; Unemployment = 30.9; UnemploymentYoy = 37.7 },
    @{ Row = 83; Province = "Prince Edward Island";     Unemployment = 14.9; UnemploymentYoy = 8.5 },
    @{ Row = 84; Province = "Nova Scotia";              Unemployment = 6.9;  UnemploymentYoy = 38.9 },
    @{ Row = 85; Province = "New Brunswick";            Unemployment = 28.9; UnemploymentYoy = 40.6 },
    @{ Row = 86; Province = "Quebec";                   Unemployment = 52.2; UnemploymentYoy = 341.9 },
    @{ Row = 87; Province = "Ontario";                  Unemployment = 80.2; UnemploymentYoy = 757.4 },
    @{ Row = 88; Province = "Manitoba";                 Unemployment = 49.8; UnemploymentYoy = 49.3 },
    @{ Row = 89; Province = "Saskatchewan";             Unemployment = 26.6; UnemploymentYoy = 41.4 },
    @{ Row = 90; Province = "Alberta";                  Unemployment = 77.7; UnemploymentYoy = 293.2 },
    @{ Row = 91; Province = "British Columbia";         Unemployment = 73.7; UnemploymentYoy = 223.7 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $wsProvince.Range("A$rowNum").Value = 44075
    $wsProvince.Range("A$rowNum").NumberFormat = "d-mmm-yy"
    $wsProvince.Range("B$rowNum").Value = $r.Province
    $wsProvince.Range("C$rowNum").Value = $r.Unemployment
    $wsProvince.Range("D$rowNum").Value = $r.UnemploymentYoy
}

# The first province row of the new date block (row 82, "Newfoundland &
# Labrador") carries the same date-style stamp on column B as every other
# first-row-of-a-date-block in the sheet.
$wsProvince.Range("B82").NumberFormat = "d-mmm-yy"

# --- View state: mirror the scrolled/selected state from the edit -------
$wsCanada.Range("C11").Select()
$wsProvince.Activate()
$wsProvince.Range("C92").Select()
